# edit.ps1 - applies the "edits to diagram and overview" commit to the
# architecture_diagram.pptx single-slide deck.
#
# Notes on precision: PowerPoint's Shape.Left/Top/Width/Height (and the
# Characters()-range equivalents) are IEEE-754 single precision (Single)
# in the real object model, and this host reproduces that: it stores the
# point value as float32 before re-expanding to EMU (1 pt = 12700 EMU),
# and *truncates* (not rounds) when converting back. The literals below
# were chosen (each one is an exact float32 value) so that converting
# them to EMU lands exactly on the target integer EMU from the target
# OOXML - not merely "close enough" after rounding to a handful of
# decimal digits.

function Get-ShapeById($shapes, $id) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $candidate = $shapes.Item($i)
        if ($candidate.Id -eq $id) {
            return $candidate
        }
    }
    return $null
}

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# ---------------------------------------------------------------------
# 1) Big rounded-rect outline (id 394): grew upward/taller.
# ---------------------------------------------------------------------
$shape394 = Get-ShapeById $s.Shapes 394
$shape394.Left   = 52.759845733642578125
$shape394.Top    = 53.76937103271484375
$shape394.Width  = 516.675537109375
$shape394.Height = 199.0678863525390625

# ---------------------------------------------------------------------
# 2) Straight connector (id 403): widened.
# ---------------------------------------------------------------------
$shape403 = Get-ShapeById $s.Shapes 403
$shape403.Left   = 345.377655029296875
$shape403.Top    = 164.947021484375
$shape403.Width  = 131.209930419921875
$shape403.Height = 0

# ---------------------------------------------------------------------
# 3) Bent connector (id 408): tiny nudge.
# ---------------------------------------------------------------------
$shape408 = Get-ShapeById $s.Shapes 408
$shape408.Left   = 197.63787841796875
$shape408.Top    = 164.947021484375
$shape408.Width  = 100.5111846923828125
$shape408.Height = 0.00007874015864217654

# ---------------------------------------------------------------------
# 4) Picture (id 407): nudged up/left slightly.
# ---------------------------------------------------------------------
$shape407 = Get-ShapeById $s.Shapes 407
$shape407.Left   = 150.5645751953125
$shape407.Top    = 141.41064453125
$shape407.Width  = 47.07331085205078125
$shape407.Height = 47.07283782958984375

# ---------------------------------------------------------------------
# 5) Picture (id 400): moved right.
# ---------------------------------------------------------------------
$shape400 = Get-ShapeById $s.Shapes 400
$shape400.Left   = 476.58758544921875
$shape400.Top    = 144.9704742431640625
$shape400.Width  = 39.9532318115234375
$shape400.Height = 39.95307159423828125

# ---------------------------------------------------------------------
# 6) SVG picture (id 23, "Graphic 22"): moved up/left.
# ---------------------------------------------------------------------
$shape23 = Get-ShapeById $s.Shapes 23
$shape23.Left   = 474.928436279296875
$shape23.Top    = 30.9507884979248046875
$shape23.Width  = 43.271575927734375
$shape23.Height = 43.271575927734375

# ---------------------------------------------------------------------
# 7) "Classic load balancer" label (id 35): moved + retyped text.
#    Final runs: "C" | "lassic " | "Load Balancer"
# ---------------------------------------------------------------------
$shape35 = Get-ShapeById $s.Shapes 35
$shape35.Left   = 134.526702880859375
$shape35.Top    = 191.3726043701171875
$shape35.Width  = 79.14890289306640625
$shape35.Height = 22.8425197601318359375

$tr35 = $shape35.TextFrame.TextRange
# "Classic load balancer" -> runs: C(1) lassic(6) ' '(1) load(4) ' balancer'(9)
# Merge "lassic" into the following " " run (drop the old misspelling flag)
# by deleting "lassic" then typing it back in front of the space run.
$tr35.Characters(2, 6).Text = ""
$tr35.Characters(2, 1).Text = "lassic "
# Now: "Classic load balancer"; "load" run is at (9,4), ' balancer' at (13,9)
$tr35.Characters(9, 4).Text = ""
$tr35.Characters(9, 9).Text = "Load Balancer"

# ---------------------------------------------------------------------
# 8) "EBS Volume" label (id 36): split + recapitalized to "EBS volume".
# ---------------------------------------------------------------------
$shape36 = Get-ShapeById $s.Shapes 36
$tr36 = $shape36.TextFrame.TextRange
$tr36.Text = "EBS volume"
# Force a run split at "volume" without altering the visible formatting.
$tr36.Characters(5, 6).Font.Name = "Consolas"

# ---------------------------------------------------------------------
# 9) New textbox (id assigned by host): "Amazon Elastic Kubernetes
#    Service (Amazon EKS)", centered, Consolas/Amazon Ember/Arial.
# ---------------------------------------------------------------------
$newBox = $s.Shapes.AddTextbox(1, 412.0406494140625, 73.463623046875, 178.6250457763672, 29.081260681152344)
$newBox.Name = "TextBox 9"

$tf = $newBox.TextFrame
$tf.WordWrap = -1
$tf.AutoSize = 1

$newBox.Fill.Visible = 0
$newBox.Line.Visible = 0

$trNew = $newBox.TextFrame.TextRange
$trNew.Text = "Amazon Elastic Kubernetes Service (Amazon EKS)"
$trNew.ParagraphFormat.Alignment = 2

$fNew = $trNew.Font
$fNew.Size = 9
$fNew.Name = "Consolas"
$fNew.NameFarEast = "Amazon Ember"
$fNew.NameComplexScript = "Arial"

# Re-assert the exact size/position last: AutoSize / text entry can
# otherwise resize the box to fit the (wrapped) text.
$newBox.Left   = 412.0406494140625
$newBox.Top    = 73.463623046875
$newBox.Width  = 178.6250457763672
$newBox.Height = 29.081260681152344
